# Generate Report for Handback
# - Marks zh-cn / de-de handback rows as "Handed back: in sync with en-US"
# - Fills in the Latest Target File (hyperlink) + Latest Handback File columns
# - Stamps the Latest Handback DateTime for each locale
# - Widens the columns that now hold longer text

$wb = $excel.ActiveWorkbook

$mdFileName = "f3621c04-1fb4-4a02-9775-7e3ee7e7defb.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d20118c3d65992a88bca8e5eb3d5335085188920/e2e/f3621c04-1fb4-4a02-9775-7e3ee7e7defb.md"
$status = "Handed back: in sync with en-US"

$zhXlf = "f3621c04-1fb4-4a02-9775-7e3ee7e7defb.4453b1d97e2c9b9d5257939f998595323278fe14.zh-cn.xlf"
$deXlf = "f3621c04-1fb4-4a02-9775-7e3ee7e7defb.4453b1d97e2c9b9d5257939f998595323278fe14.de-de.xlf"

$zhHandbackDate = "2016-09-04 01:04:41"
$deHandbackDate = "2016-09-04 01:04:48"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, $null, $mdFileName, $mdFileName) | Out-Null
$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = $zhHandbackDate

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, $null, $mdFileName, $mdFileName) | Out-Null
$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = $deHandbackDate

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---- Overview sheet (summary columns for zh-cn / de-de) ----
# The overview's per-locale status cells mirror the same "Ready for
# handoff" text that has now been superseded everywhere, so refresh them
# to match and widen the columns that now hold the longer text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668
